$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" date column (C) for rows 2-11 from 45207 to 45208
for ($r = 2; $r -le 11; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45207) {
        $cell.Value2 = 45208
    }
}

# Update hyperlink formulas in row 2 (columns S, T, V, W, X, Y):
# replace "Logging_HOFORS" with "Logging_2104" in the folder paths
$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2104/artfynd/A 33491-2023.xlsx", "A 33491-2023")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2104/kartor/A 33491-2023.png", "A 33491-2023")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2104/klagomål/A 33491-2023.docx", "A 33491-2023")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2104/klagomålsmail/A 33491-2023.docx", "A 33491-2023")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2104/tillsyn/A 33491-2023.docx", "A 33491-2023")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2104/tillsynsmail/A 33491-2023.docx", "A 33491-2023")'
